$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-12 down to 10-13
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the weekly price record
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44449
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = "Bruselas (repollito)"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 220
$ws.Range("K9").Value = 22000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 23091
$ws.Range("N9").Value = "$/malla 15 kilos"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 1539
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"
